# Updates the cryptos list values (prices / 1h volume %) and fixes the
# Kaspa / Fetch.AI row ordering, matching the latest data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force Excel to store the value as text (not auto-convert to a
    # number) by using a leading apostrophe, then reset the cell style
    # back to Normal so no stray number-format style is introduced.
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.807.76"
$ws.Range("E2").Value = "  +0.81%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.678.77"
$ws.Range("E3").Value = "  +0.89%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
Set-TextValue "D5" "601.81"
$ws.Range("E5").Value = "  -0.59%  "

# Row 6 - Solana
Set-TextValue "D6" "156.48"
$ws.Range("E6").Value = "  -0.90%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.07%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.22%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +0.18%  "

# Row 10 - Toncoin
Set-TextValue "D10" "5.91"
$ws.Range("E10").Value = "  +1.55%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.398"
$ws.Range("E11").Value = "  -3.63%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.17%  "

# Row 13 - Avalanche
Set-TextValue "D13" "29.48"
$ws.Range("E13").Value = "  -0.40%  "

# Row 14 - ShibaInu
Set-TextValue "D14" "0.0000203"
$ws.Range("E14").Value = "  +6.64%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.156.94"
$ws.Range("E15").Value = "  +0.65%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "65.605.99"
$ws.Range("E16").Value = "  +0.75%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.686.44"
$ws.Range("E17").Value = "  +1.10%  "

# Row 18 - Chainlink
Set-TextValue "D18" "12.61"
$ws.Range("E18").Value = "  -1.14%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -1.99%  "

# Row 20 - Uniswap
Set-TextValue "D20" "7.60"
$ws.Range("E20").Value = "  +2.96%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "351.44"
$ws.Range("E21").Value = "  -2.38%  "

# Row 22 - Dai
Set-TextValue "D22" "0.999"
$ws.Range("E22").Value = "  -0.09%  "

# Row 23 - Litecoin
Set-TextValue "D23" "70.06"
$ws.Range("E23").Value = "  +1.02%  "

# Row 24 - PEPE
Set-TextValue "D24" "0.0000111"
$ws.Range("E24").Value = "  +6.75%  "

# Row 25 - InternetComputer(DFINITY)
Set-TextValue "D25" "9.84"
$ws.Range("E25").Value = "  +2.37%  "

# Row 26 - SuiNetwork
Set-TextValue "D26" "1.63"
$ws.Range("E26").Value = "  -5.58%  "

# Row 27 - was Kaspa, now Fetch.AI
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D27" "1.63"
$ws.Range("E27").Value = "  -2.06%  "

# Row 28 - was Fetch.AI, now Kaspa
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D28" "0.169"
$ws.Range("E28").Value = "  +2.36%  "

# Row 29 - Aptos
Set-TextValue "D29" "8.17"
$ws.Range("E29").Value = "  -1.38%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  +0.01%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -2.85%  "

# Row 32 - Bittensor
Set-TextValue "D32" "529.55"
$ws.Range("E32").Value = "  -3.59%  "

# Row 33 - ImmutableX
$ws.Range("E33").Value = "  -4.00%  "

# Row 34 - RenderToken
Set-TextValue "D34" "6.57"
$ws.Range("E34").Value = "  +2.11%  "

# Row 35 - NEARProtocol
Set-TextValue "D35" "5.42"
$ws.Range("E35").Value = "  -4.23%  "

# Row 36 - PolygonEcosystemToken
$ws.Range("E36").Value = "  -1.68%  "

# Row 37 - EthereumClassic
Set-TextValue "D37" "20.48"
$ws.Range("E37").Value = "  -0.62%  "

# Row 38 - Monero
Set-TextValue "D38" "160.90"
$ws.Range("E38").Value = "  -1.72%  "

# Row 39 - FirstDigitalUSD
Set-TextValue "D39" "0.999"
$ws.Range("E39").Value = "  -0.02%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  -2.30%  "

# Row 41 - USDe
Set-TextValue "D41" "0.999"
$ws.Range("E41").Value = "  +0.03%  "

# Row 42 - OKB
Set-TextValue "D42" "42.25"
$ws.Range("E42").Value = "  -0.62%  "

# Row 43 - Aave
Set-TextValue "D43" "165.98"
$ws.Range("E43").Value = "  -1.17%  "

# Row 44 - Filecoin
Set-TextValue "D44" "4.10"
$ws.Range("E44").Value = "  -2.62%  "

# Row 45 - Hedera
Set-TextValue "D45" "0.0623"
$ws.Range("E45").Value = "  +0.12%  "

# Row 46 - InjectiveProtocol
Set-TextValue "D46" "23.13"
$ws.Range("E46").Value = "  -0.35%  "

# Row 47 - dogwifhat
$ws.Range("E47").Value = "  -4.06%  "

# Row 48 - VeChain
$ws.Range("E48").Value = "  -0.51%  "

# Row 49 - Mantle
$ws.Range("E49").Value = "  -1.49%  "

# Row 50 - EnergySwap
Set-TextValue "D50" "20.30"
$ws.Range("E50").Value = "  +2.27%  "

# Row 51 - Stellar
Set-TextValue "D51" "0.0986"
$ws.Range("E51").Value = "  -0.01%  "
